$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data rows to append (datetime serial, open/high/low/close value)
$newRows = @(
    @{ Row = 291; Date = 44926; Value = 41678378000000 },
    @{ Row = 292; Date = 44957; Value = 42206471000000 },
    @{ Row = 293; Date = 44985; Value = 41767936000000 },
    @{ Row = 294; Date = 45016; Value = 42100727000000 },
    @{ Row = 295; Date = 45046; Value = 41265420000000 },
    @{ Row = 296; Date = 45077; Value = 40642084000000 },
    @{ Row = 297; Date = 45107; Value = 41806284000000 }
)

# Use the last existing data row (290) as a formatting template so the
# appended rows inherit the same cell styles already used in the sheet.
$templateRow = 290

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$templateRow`:G$templateRow").Copy()
    $ws.Range("A$row`:G$row").PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = "ECONOMICS:CNCBBS"
    $ws.Cells.Item($row, 3).Value = $r.Value
    $ws.Cells.Item($row, 4).Value = $r.Value
    $ws.Cells.Item($row, 5).Value = $r.Value
    $ws.Cells.Item($row, 6).Value = $r.Value
    $ws.Cells.Item($row, 7).Value = 0
}

$excel.CutCopyMode = 0
